$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Step 1: insert the new "accuracy_balanced_mean" block (7 rows) before the old f1_macro_std block ----
$ws.Range("A16:A22").EntireRow.Insert()

# Copy the row formatting (bold/border/center-aligned col A, plain B-H) from an existing 7-row block onto the new rows
$ws.Range("A2:H8").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- Step 2: append the new "accuracy_balanced_std" block (7 rows) after the (now shifted) f1_micro_std block ----
$ws.Range("A37:A43").EntireRow.Insert()
$ws.Range("A2:H8").Copy()
$ws.Range("A37").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- Step 3: write cell values for the two new blocks ----
$ws.Range("A16").Value = "accuracy_balanced_mean"
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "0"
$ws.Range("C16").Value = 0
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 0.132

$ws.Range("A17").Value = "accuracy_balanced_mean"
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "100"
$ws.Range("C17").Value = 0.174
$ws.Range("D17").Value = 0.142
$ws.Range("E17").Value = 0.277
$ws.Range("F17").Value = 0.291
$ws.Range("G17").Value = 0.236
$ws.Range("H17").Value = 0.338

$ws.Range("A18").Value = "accuracy_balanced_mean"
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "500"
$ws.Range("C18").Value = 0.28
$ws.Range("D18").Value = 0.301
$ws.Range("E18").Value = 0.367
$ws.Range("F18").Value = 0.381
$ws.Range("G18").Value = 0.281
$ws.Range("H18").Value = 0.467

$ws.Range("A19").Value = "accuracy_balanced_mean"
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = "1000"
$ws.Range("C19").Value = 0.31
$ws.Range("D19").Value = 0.323
$ws.Range("E19").Value = 0.404
$ws.Range("F19").Value = 0.409
$ws.Range("G19").Value = 0.421
$ws.Range("H19").Value = 0.486

$ws.Range("A20").Value = "accuracy_balanced_mean"
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "2500"
$ws.Range("C20").Value = 0.347
$ws.Range("D20").Value = 0.359
$ws.Range("E20").Value = 0.429
$ws.Range("F20").Value = 0.435
$ws.Range("G20").Value = 0.473
$ws.Range("H20").Value = 0.5

$ws.Range("A21").Value = "accuracy_balanced_mean"
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "5000"
$ws.Range("C21").Value = 0.376
$ws.Range("D21").Value = 0.371
$ws.Range("E21").Value = 0.445
$ws.Range("F21").Value = 0.451
$ws.Range("G21").Value = 0.503
$ws.Range("H21").Value = 0.534

$ws.Range("A22").Value = "accuracy_balanced_mean"
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = "10000"
$ws.Range("C22").Value = 0.397
$ws.Range("D22").Value = 0.404
$ws.Range("E22").Value = 0.46
$ws.Range("F22").Value = 0.464
$ws.Range("G22").Value = 0.517
$ws.Range("H22").Value = 0.548

$ws.Range("A37").Value = "accuracy_balanced_std"
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "0"
$ws.Range("C37").Value = 0
$ws.Range("D37").Value = 0
$ws.Range("E37").Value = 0
$ws.Range("F37").Value = 0
$ws.Range("G37").Value = 0
$ws.Range("H37").Value = 0

$ws.Range("A38").Value = "accuracy_balanced_std"
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "100"
$ws.Range("C38").Value = 0.003
$ws.Range("D38").Value = 0.004
$ws.Range("E38").Value = 0.002
$ws.Range("F38").Value = 0.007
$ws.Range("G38").Value = 0.021
$ws.Range("H38").Value = 0.012

$ws.Range("A39").Value = "accuracy_balanced_std"
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "500"
$ws.Range("C39").Value = 0.01
$ws.Range("D39").Value = 0.01
$ws.Range("E39").Value = 0.009
$ws.Range("F39").Value = 0.006
$ws.Range("G39").Value = 0.111
$ws.Range("H39").Value = 0.013

$ws.Range("A40").Value = "accuracy_balanced_std"
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "1000"
$ws.Range("C40").Value = 0.011
$ws.Range("D40").Value = 0.012
$ws.Range("E40").Value = 0.007
$ws.Range("F40").Value = 0.009
$ws.Range("G40").Value = 0.016
$ws.Range("H40").Value = 0.007

$ws.Range("A41").Value = "accuracy_balanced_std"
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "2500"
$ws.Range("C41").Value = 0.004
$ws.Range("D41").Value = 0.002
$ws.Range("E41").Value = 0.001
$ws.Range("F41").Value = 0.006
$ws.Range("G41").Value = 0.006
$ws.Range("H41").Value = 0.007

$ws.Range("A42").Value = "accuracy_balanced_std"
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "5000"
$ws.Range("C42").Value = 0.004
$ws.Range("D42").Value = 0.004
$ws.Range("E42").Value = 0.002
$ws.Range("F42").Value = 0.006
$ws.Range("G42").Value = 0.004
$ws.Range("H42").Value = 0.005

$ws.Range("A43").Value = "accuracy_balanced_std"
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "10000"
$ws.Range("C43").Value = 0.001
$ws.Range("D43").Value = 0.002
$ws.Range("E43").Value = 0.002
$ws.Range("F43").Value = 0.005
$ws.Range("G43").Value = 0.013
$ws.Range("H43").Value = 0.004

# ---- Step 4: normalize column B style back to default (no explicit style) for the new rows ----
$ws.Range("B2").Copy()
$ws.Range("B16:B22").PasteSpecial(-4122)
$ws.Range("B2").Copy()
$ws.Range("B37:B43").PasteSpecial(-4122)
$excel.CutCopyMode = $false

Write-Host "done"